# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the 570d0173-...md file row
# (row 6) in both the zh-cn and de-de localization-status sheets, as a
# fresh handoff report is generated (new handoff timestamps recorded).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-24 14:51:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-24 14:52:04"
